$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Update existing row 110 (G110 24->22, H110 48->60). K110/L110 are
#    formulas and will recompute automatically.
# ---------------------------------------------------------------------------
$ws.Cells.Item(110, 7).Value = 22
$ws.Cells.Item(110, 8).Value = 60

# ---------------------------------------------------------------------------
# 2. Add two new rows (111: Y10 "Super Sentinel", 112: Drewry Shunter)
# ---------------------------------------------------------------------------

# --- Row 111 ---------------------------------------------------------------
$ws.Cells.Item(111, 1).Value = "Y10 ""Super Sentinel"""
$ws.Cells.Item(111, 2).Value = 1930
$ws.Cells.Item(111, 3).Formula = "=B111-B110"
$ws.Cells.Item(111, 4).Value = 1
$ws.Cells.Item(111, 5).Value = "Cargo Tram"
$ws.Cells.Item(111, 6).Value = 16041
$ws.Cells.Item(111, 7).Value = 24
$ws.Cells.Item(111, 8).Value = 60
$ws.Cells.Item(111, 9).Value = 26
$ws.Cells.Item(111, 10).Value = 200

$ws.Cells.Item(111, 11).NumberFormat = "0"
$ws.Cells.Item(111, 11).Formula = "=SQRT(G111*H111)*POWER((MIN(I111,20)+SQRT(MAX(I111-20,0))),0.9)*`$B`$1"

$ws.Cells.Item(111, 12).NumberFormat = "0"
$ws.Cells.Item(111, 12).Formula = "=POWER((G111*G111*H111), 0.33)*LOG10(J111)*10*`$B`$1"

$ws.Cells.Item(111, 13).NumberFormat = "0"
$ws.Cells.Item(111, 13).Value = "x"

$ws.Cells.Item(111, 14).Value = 7

$ws.Cells.Item(111, 16).NumberFormat = "0.0"

$ws.Cells.Item(111, 17).NumberFormat = "0"
$ws.Cells.Item(111, 17).Formula = "=CONCATENATE(ROUND(N111*VLOOKUP(E111,'ID Scheme'!`$A`$2:`$E`$7,3),0), ""x"",ROUND(O111*VLOOKUP(E111,'ID Scheme'!`$A`$2:`$E`$7,5),0), ""x"",ROUND(P111*VLOOKUP(E111,'ID Scheme'!`$A`$2:`$E`$7,4),0))"

# --- Row 112 -----------------------------------------------------------------
$ws.Cells.Item(112, 1).Value = "Drewry Shunter"
$ws.Cells.Item(112, 2).Value = 1952
$ws.Cells.Item(112, 3).Formula = "=B112-B111"
$ws.Cells.Item(112, 4).Value = 2
$ws.Cells.Item(112, 5).Value = "Cargo Tram"
$ws.Cells.Item(112, 6).Value = 16051
$ws.Cells.Item(112, 7).Value = 27
$ws.Cells.Item(112, 8).Value = 60
$ws.Cells.Item(112, 9).Value = 26
$ws.Cells.Item(112, 10).Value = 152

$ws.Cells.Item(112, 11).NumberFormat = "0"
$ws.Cells.Item(112, 11).Formula = "=SQRT(G112*H112)*POWER((MIN(I112,20)+SQRT(MAX(I112-20,0))),0.9)*`$B`$1"

$ws.Cells.Item(112, 12).NumberFormat = "0"
$ws.Cells.Item(112, 12).Formula = "=POWER((G112*G112*H112), 0.33)*LOG10(J112)*10*`$B`$1"

$ws.Cells.Item(112, 13).NumberFormat = "0"
$ws.Cells.Item(112, 13).Value = "x"

$ws.Cells.Item(112, 14).Value = 7

$ws.Cells.Item(112, 16).NumberFormat = "0.0"

$ws.Cells.Item(112, 17).NumberFormat = "0"
$ws.Cells.Item(112, 17).Formula = "=CONCATENATE(ROUND(N112*VLOOKUP(E112,'ID Scheme'!`$A`$2:`$E`$7,3),0), ""x"",ROUND(O112*VLOOKUP(E112,'ID Scheme'!`$A`$2:`$E`$7,5),0), ""x"",ROUND(P112*VLOOKUP(E112,'ID Scheme'!`$A`$2:`$E`$7,4),0))"

# ---------------------------------------------------------------------------
# 3. Conditional formatting: the rule on column C (highlight gaps > 10) is
#    split so rows 111 and 112 get their own copy of the "red" rule.
# ---------------------------------------------------------------------------
$r111 = $ws.Range("C111")
$fc111 = $r111.FormatConditions.Add(1, 5, "10")
$fc111.Font.Color = 393372
$fc111.Interior.Color = 13551615

$r112 = $ws.Range("C112")
$fc112 = $r112.FormatConditions.Add(1, 5, "10")
$fc112.Font.Color = 393372
$fc112.Interior.Color = 13551615

$fcs0 = $ws.Range("C1").FormatConditions
$fc0 = $fcs0.Item(1)
$fc0.Priority = 10
$fc111.Priority = 2
$fc112.Priority = 1

# ---------------------------------------------------------------------------
# 4. View state: selection moves to F111, view scrolled so row 85 is the
#    first visible row below the frozen header rows (1-3 stay frozen).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F111").Select()
